# This script updates the NATMI ligand-receptor summary (Hgf -> St14)
# with recalculated values reflecting an updated TPM expression matrix.
# Only the Ligand/Receptor average & total expression values (columns G, H, M, N)
# changed for the clusters whose TPM was revised (Inflammatory-Mac, MuSCs,
# Resolving-Mac on the ligand side; ECs, Inflammatory-Mac, MuSCs, Resolving-Mac
# on the receptor side). All specificity and edge-weight columns
# (I, J, O, P, Q, R, S, T) are derivative of those and are updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.001628901849080777
$ws.Range("J2").Value = 0.001628901849080777
$ws.Range("M2").Value = 1.424886333333333
$ws.Range("N2").Value = 4.274659
$ws.Range("O2").Value = 0.3775790263072122
$ws.Range("P2").Value = 0.3775790263072122
$ws.Range("Q2").Value = 0.1334344306092222
$ws.Range("R2").Value = 1.200909875483
$ws.Range("S2").Value = 0.0006150391741259374
$ws.Range("T2").Value = 0.0006150391741259373

# Row 3
$ws.Range("I3").Value = 0.001628901849080777
$ws.Range("J3").Value = 0.001628901849080777
$ws.Range("O3").Value = 0.09562220712767076
$ws.Range("P3").Value = 0.09562220712767076
$ws.Range("S3").Value = 0.000155759190003448
$ws.Range("T3").Value = 0.0001557591900034479

# Row 4
$ws.Range("I4").Value = 0.001628901849080777
$ws.Range("J4").Value = 0.001628901849080777
$ws.Range("M4").Value = 0.6353876666666666
$ws.Range("N4").Value = 1.906163
$ws.Range("O4").Value = 0.1683706629050024
$ws.Range("P4").Value = 0.1683706629050024
$ws.Range("Q4").Value = 0.05950130163677777
$ws.Range("R4").Value = 0.5355117147309999
$ws.Range("S4").Value = 0.0002742592841369145
$ws.Range("T4").Value = 0.0002742592841369145

# Row 5
$ws.Range("I5").Value = 0.001628901849080777
$ws.Range("J5").Value = 0.001628901849080777
$ws.Range("M5").Value = 0.4155976666666667
$ws.Range("N5").Value = 1.246793
$ws.Range("O5").Value = 0.110128758094306
$ws.Range("P5").Value = 0.110128758094306
$ws.Range("Q5").Value = 0.03891892056011111
$ws.Range("R5").Value = 0.350270285041
$ws.Range("S5").Value = 0.0001793889376967847
$ws.Range("T5").Value = 0.0001793889376967846

# Row 6
$ws.Range("I6").Value = 0.001628901849080777
$ws.Range("J6").Value = 0.001628901849080777
$ws.Range("M6").Value = 0.937018
$ws.Range("N6").Value = 2.811054
$ws.Range("O6").Value = 0.2482993455658087
$ws.Range("P6").Value = 0.2482993455658087
$ws.Range("Q6").Value = 0.08774767528866667
$ws.Range("R6").Value = 0.7897290775979999
$ws.Range("S6").Value = 0.0004044552631176926
$ws.Range("T6").Value = 0.0004044552631176925

# Row 7
$ws.Range("I7").Value = 0.1785947081647151
$ws.Range("J7").Value = 0.178594708164715
$ws.Range("M7").Value = 1.424886333333333
$ws.Range("N7").Value = 4.274659
$ws.Range("O7").Value = 0.3775790263072122
$ws.Range("P7").Value = 0.3775790263072122
$ws.Range("Q7").Value = 14.62990738651756
$ws.Range("R7").Value = 131.669166478658
$ws.Range("S7").Value = 0.06743361601245383
$ws.Range("T7").Value = 0.06743361601245382

# Row 8
$ws.Range("I8").Value = 0.1785947081647151
$ws.Range("J8").Value = 0.178594708164715
$ws.Range("O8").Value = 0.09562220712767076
$ws.Range("P8").Value = 0.09562220712767076
$ws.Range("S8").Value = 0.0170776201760323
$ws.Range("T8").Value = 0.01707762017603229

# Row 9
$ws.Range("I9").Value = 0.1785947081647151
$ws.Range("J9").Value = 0.178594708164715
$ws.Range("M9").Value = 0.6353876666666666
$ws.Range("N9").Value = 1.906163
$ws.Range("O9").Value = 0.1683706629050024
$ws.Range("P9").Value = 0.1683706629050024
$ws.Range("Q9").Value = 6.523792460078443
$ws.Range("R9").Value = 58.71413214070599
$ws.Range("S9").Value = 0.03007010940501851
$ws.Range("T9").Value = 0.03007010940501851

# Row 10
$ws.Range("I10").Value = 0.1785947081647151
$ws.Range("J10").Value = 0.178594708164715
$ws.Range("M10").Value = 0.4155976666666667
$ws.Range("N10").Value = 1.246793
$ws.Range("O10").Value = 0.110128758094306
$ws.Range("P10").Value = 0.110128758094306
$ws.Range("Q10").Value = 4.267116071751778
$ws.Range("R10").Value = 38.404044645766
$ws.Range("S10").Value = 0.01966841341239509
$ws.Range("T10").Value = 0.01966841341239508

# Row 11
$ws.Range("I11").Value = 0.1785947081647151
$ws.Range("J11").Value = 0.178594708164715
$ws.Range("M11").Value = 0.937018
$ws.Range("N11").Value = 2.811054
$ws.Range("O11").Value = 0.2482993455658087
$ws.Range("P11").Value = 0.2482993455658087
$ws.Range("Q11").Value = 9.620757978238666
$ws.Range("R11").Value = 86.586821804148
$ws.Range("S11").Value = 0.04434494915881534
$ws.Range("T11").Value = 0.04434494915881533

# Row 12
$ws.Range("G12").Value = 23.67539566666666
$ws.Range("H12").Value = 71.02618699999999
$ws.Range("I12").Value = 0.4118171950916292
$ws.Range("J12").Value = 0.4118171950916292
$ws.Range("M12").Value = 1.424886333333333
$ws.Range("N12").Value = 4.274659
$ws.Range("O12").Value = 0.3775790263072122
$ws.Range("P12").Value = 0.3775790263072122
$ws.Range("Q12").Value = 33.73474772169255
$ws.Range("R12").Value = 303.6127294952329
$ws.Range("S12").Value = 0.1554935355392646
$ws.Range("T12").Value = 0.1554935355392646

# Row 13
$ws.Range("G13").Value = 23.67539566666666
$ws.Range("H13").Value = 71.02618699999999
$ws.Range("I13").Value = 0.4118171950916292
$ws.Range("J13").Value = 0.4118171950916292
$ws.Range("O13").Value = 0.09562220712767076
$ws.Range("P13").Value = 0.09562220712767076
$ws.Range("Q13").Value = 8.543353336100777
$ws.Range("R13").Value = 76.890180024907
$ws.Range("S13").Value = 0.03937886912778817
$ws.Range("T13").Value = 0.03937886912778817

# Row 14
$ws.Range("G14").Value = 23.67539566666666
$ws.Range("H14").Value = 71.02618699999999
$ws.Range("I14").Value = 0.4118171950916292
$ws.Range("J14").Value = 0.4118171950916292
$ws.Range("M14").Value = 0.6353876666666666
$ws.Range("N14").Value = 1.906163
$ws.Range("O14").Value = 0.1683706629050024
$ws.Range("P14").Value = 0.1683706629050024
$ws.Range("Q14").Value = 15.04305441005344
$ws.Range("R14").Value = 135.387489690481
$ws.Range("S14").Value = 0.0693379341332563
$ws.Range("T14").Value = 0.06933793413325628

# Row 15
$ws.Range("G15").Value = 23.67539566666666
$ws.Range("H15").Value = 71.02618699999999
$ws.Range("I15").Value = 0.4118171950916292
$ws.Range("J15").Value = 0.4118171950916292
$ws.Range("M15").Value = 0.4155976666666667
$ws.Range("N15").Value = 1.246793
$ws.Range("O15").Value = 0.110128758094306
$ws.Range("P15").Value = 0.110128758094306
$ws.Range("Q15").Value = 9.839439196476777
$ws.Range("R15").Value = 88.55495276829099
$ws.Range("S15").Value = 0.04535291625732166
$ws.Range("T15").Value = 0.04535291625732166

# Row 16
$ws.Range("G16").Value = 23.67539566666666
$ws.Range("H16").Value = 71.02618699999999
$ws.Range("I16").Value = 0.4118171950916292
$ws.Range("J16").Value = 0.4118171950916292
$ws.Range("M16").Value = 0.937018
$ws.Range("N16").Value = 2.811054
$ws.Range("O16").Value = 0.2482993455658087
$ws.Range("P16").Value = 0.2482993455658087
$ws.Range("Q16").Value = 22.18427189678867
$ws.Range("R16").Value = 199.658447071098
$ws.Range("S16").Value = 0.1022539400339985
$ws.Range("T16").Value = 0.1022539400339985

# Row 17
$ws.Range("G17").Value = 0.3314846666666666
$ws.Range("H17").Value = 0.9944539999999999
$ws.Range("I17").Value = 0.005765947381177186
$ws.Range("J17").Value = 0.005765947381177185
$ws.Range("M17").Value = 1.424886333333333
$ws.Range("N17").Value = 4.274659
$ws.Range("O17").Value = 0.3775790263072122
$ws.Range("P17").Value = 0.3775790263072122
$ws.Range("Q17").Value = 0.4723279712428888
$ws.Range("R17").Value = 4.250951741185999
$ws.Range("S17").Value = 0.002177100797923502
$ws.Range("T17").Value = 0.002177100797923502

# Row 18
$ws.Range("G18").Value = 0.3314846666666666
$ws.Range("H18").Value = 0.9944539999999999
$ws.Range("I18").Value = 0.005765947381177186
$ws.Range("J18").Value = 0.005765947381177185
$ws.Range("O18").Value = 0.09562220712767076
$ws.Range("P18").Value = 0.09562220712767076
$ws.Range("Q18").Value = 0.1196174574104444
$ws.Range("R18").Value = 1.076557116694
$ws.Range("S18").Value = 0.0005513526147701757
$ws.Range("T18").Value = 0.0005513526147701756

# Row 19
$ws.Range("G19").Value = 0.3314846666666666
$ws.Range("H19").Value = 0.9944539999999999
$ws.Range("I19").Value = 0.005765947381177186
$ws.Range("J19").Value = 0.005765947381177185
$ws.Range("M19").Value = 0.6353876666666666
$ws.Range("N19").Value = 1.906163
$ws.Range("O19").Value = 0.1683706629050024
$ws.Range("P19").Value = 0.1683706629050024
$ws.Range("Q19").Value = 0.2106212688891111
$ws.Range("R19").Value = 1.895591420002
$ws.Range("S19").Value = 0.0009708163828441651
$ws.Range("T19").Value = 0.0009708163828441649

# Row 20
$ws.Range("G20").Value = 0.3314846666666666
$ws.Range("H20").Value = 0.9944539999999999
$ws.Range("I20").Value = 0.005765947381177186
$ws.Range("J20").Value = 0.005765947381177185
$ws.Range("M20").Value = 0.4155976666666667
$ws.Range("N20").Value = 1.246793
$ws.Range("O20").Value = 0.110128758094306
$ws.Range("P20").Value = 0.110128758094306
$ws.Range("Q20").Value = 0.1377642540024444
$ws.Range("R20").Value = 1.239878286022
$ws.Range("S20").Value = 0.0006349966243261595
$ws.Range("T20").Value = 0.0006349966243261594

# Row 21
$ws.Range("G21").Value = 0.3314846666666666
$ws.Range("H21").Value = 0.9944539999999999
$ws.Range("I21").Value = 0.005765947381177186
$ws.Range("J21").Value = 0.005765947381177185
$ws.Range("M21").Value = 0.937018
$ws.Range("N21").Value = 2.811054
$ws.Range("O21").Value = 0.2482993455658087
$ws.Range("P21").Value = 0.2482993455658087
$ws.Range("Q21").Value = 0.3106070993906667
$ws.Range("R21").Value = 2.795463894516
$ws.Range("S21").Value = 0.001431680961313183
$ws.Range("T21").Value = 0.001431680961313183

# Row 22
$ws.Range("G22").Value = 23.12211433333333
$ws.Range("H22").Value = 69.366343
$ws.Range("I22").Value = 0.4021932475133977
$ws.Range("J22").Value = 0.4021932475133977
$ws.Range("M22").Value = 1.424886333333333
$ws.Range("N22").Value = 4.274659
$ws.Range("O22").Value = 0.3775790263072122
$ws.Range("P22").Value = 0.3775790263072122
$ws.Range("Q22").Value = 32.94638471133744
$ws.Range("R22").Value = 296.517462402037
$ws.Range("S22").Value = 0.1518597347834443
$ws.Range("T22").Value = 0.1518597347834443

# Row 23
$ws.Range("G23").Value = 23.12211433333333
$ws.Range("H23").Value = 69.366343
$ws.Range("I23").Value = 0.4021932475133977
$ws.Range("J23").Value = 0.4021932475133977
$ws.Range("O23").Value = 0.09562220712767076
$ws.Range("P23").Value = 0.09562220712767076
$ws.Range("Q23").Value = 8.343699738269223
$ws.Range("R23").Value = 75.093297644423
$ws.Range("S23").Value = 0.03845860601907667
$ws.Range("T23").Value = 0.03845860601907667

# Row 24
$ws.Range("G24").Value = 23.12211433333333
$ws.Range("H24").Value = 69.366343
$ws.Range("I24").Value = 0.4021932475133977
$ws.Range("J24").Value = 0.4021932475133977
$ws.Range("M24").Value = 0.6353876666666666
$ws.Range("N24").Value = 1.906163
$ws.Range("O24").Value = 0.1683706629050024
$ws.Range("P24").Value = 0.1683706629050024
$ws.Range("Q24").Value = 14.69150627465655
$ws.Range("R24").Value = 132.223556471909
$ws.Range("S24").Value = 0.06771754369974647
$ws.Range("T24").Value = 0.06771754369974646

# Row 25
$ws.Range("G25").Value = 23.12211433333333
$ws.Range("H25").Value = 69.366343
$ws.Range("I25").Value = 0.4021932475133977
$ws.Range("J25").Value = 0.4021932475133977
$ws.Range("M25").Value = 0.4155976666666667
$ws.Range("N25").Value = 1.246793
$ws.Range("O25").Value = 0.110128758094306
$ws.Range("P25").Value = 0.110128758094306
$ws.Range("Q25").Value = 9.609496765333223
$ws.Range("R25").Value = 86.485470887999
$ws.Range("S25").Value = 0.04429304286256632
$ws.Range("T25").Value = 0.04429304286256631

# Row 26
$ws.Range("G26").Value = 23.12211433333333
$ws.Range("H26").Value = 69.366343
$ws.Range("I26").Value = 0.4021932475133977
$ws.Range("J26").Value = 0.4021932475133977
$ws.Range("M26").Value = 0.937018
$ws.Range("N26").Value = 2.811054
$ws.Range("O26").Value = 0.2482993455658087
$ws.Range("P26").Value = 0.2482993455658087
$ws.Range("Q26").Value = 21.66583732839133
$ws.Range("R26").Value = 194.992535955522
$ws.Range("S26").Value = 0.09986432014856396
$ws.Range("T26").Value = 0.09986432014856395
